$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 28, shifting existing rows 28..140 down to 29..141.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28. Columns A,B,C,E,F,G,H,I,J,K,L,R carry the
# same values as the row directly below (which used to be row 28 before the
# insert), while D,M,N,O,P,Q,S,T take the new weekly observation's values.
$ws.Cells.Item(28, 1).Value = $ws.Cells.Item(29, 1).Value()
$ws.Cells.Item(28, 2).Value = $ws.Cells.Item(29, 2).Value()
$ws.Cells.Item(28, 3).Value = $ws.Cells.Item(29, 3).Value()
$ws.Cells.Item(28, 4).Value = 45250
$ws.Cells.Item(28, 5).Value = $ws.Cells.Item(29, 5).Value()
$ws.Cells.Item(28, 6).Value = $ws.Cells.Item(29, 6).Value()
$ws.Cells.Item(28, 7).Value = $ws.Cells.Item(29, 7).Value()
$ws.Cells.Item(28, 8).Value = $ws.Cells.Item(29, 8).Value()
$ws.Cells.Item(28, 9).Value = $ws.Cells.Item(29, 9).Value()
$ws.Cells.Item(28, 10).Value = $ws.Cells.Item(29, 10).Value()
$ws.Cells.Item(28, 11).Value = $ws.Cells.Item(29, 11).Value()
$ws.Cells.Item(28, 12).Value = $ws.Cells.Item(29, 12).Value()
$ws.Cells.Item(28, 13).Value = 65
$ws.Cells.Item(28, 14).Value = 2600
$ws.Cells.Item(28, 15).Value = 2600
$ws.Cells.Item(28, 16).Value = 2600
$ws.Cells.Item(28, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(28, 18).Value = $ws.Cells.Item(29, 18).Value()
$ws.Cells.Item(28, 19).Value = 2600
$ws.Cells.Item(28, 20).Value = 1
